$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted at row 446, pushing the previously
# existing rows 446-553 down to 447-554 (the dataset keeps growing with
# one additional week of price data each time).
$ws.Rows.Item(446).Insert()

$ws.Range("A446").Value = 6
$ws.Range("B446").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C446").Value = "Metropolitana"
$ws.Range("D446").Value = 44943
$ws.Range("E446").Value = 13
$ws.Range("F446").Value = 100112043
$ws.Range("G446").Value = "Pepino ensalada"
$ws.Range("H446").Value = "Sin especificar"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 480
$ws.Range("K446").Value = 9000
$ws.Range("L446").Value = 9000
$ws.Range("M446").Value = 9000
$ws.Range("N446").Value = "`$/caja 60 unidades"
$ws.Range("O446").Value = "Región Metropolitana"
$ws.Range("P446").Value = 150
$ws.Range("Q446").Value = 60
$ws.Range("R446").Value = "Hortaliza"
